# Add new power plants to Electricity Source subscript (issues #280 and #99)

$wb = $excel.ActiveWorkbook

# --- "BGCL" sheet: append new rows for additional power plant types ---
$wsBgcl = $wb.Worksheets.Item("BGCL")
$wsBgcl.Select()

$wsBgcl.Range("A19").Value = "hard coal w CCS"
$wsBgcl.Range("B19").Formula = "=B2"

$wsBgcl.Range("A20").Value = "natural gas combined cycle w CCS"
$wsBgcl.Range("B20").Formula = "=B4"

$wsBgcl.Range("A21").Value = "biomass w CCS"
$wsBgcl.Range("B21").Formula = "=B10"

$wsBgcl.Range("A22").Value = "lignite w CCS"
$wsBgcl.Range("B22").Formula = "=B14"

$wsBgcl.Range("A23").Value = "small modular reactor"
$wsBgcl.Range("B23").Formula = "=B5"

$wsBgcl.Range("A24").Value = "hydrogen"
$wsBgcl.Range("B24").Formula = "=B4"

$wsBgcl.Range("A25").Select()

# --- "About" sheet: add hyperlink on B7 (already contains the EIA source URL text) ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Select()

$url = $wsAbout.Range("B7").Value2

# stash B7's current formatting so we can restore it after Hyperlinks.Add
# (which otherwise reapplies its own "Hyperlink" style variant)
$wsAbout.Range("D100").Value = "tmp"
$wsAbout.Range("B7").Copy()
$wsAbout.Range("D100").PasteSpecial(-4122)  # xlPasteFormats

$wsAbout.Hyperlinks.Add($wsAbout.Range("B7"), $url)

$wsAbout.Range("D100").Copy()
$wsAbout.Range("B7").PasteSpecial(-4122)    # xlPasteFormats
$wsAbout.Range("D100").Clear()

# remove the bold "applyFont" style override on A11 (falls back to default style)
$wsAbout.Range("A11").Style = "Normal"

$wsAbout.Range("B7").Select()
